{"js": "// Remove the placeholder \"AuthorNames\" and \"Affiliations\" paragraphs from the\n// document header block, and drop the now-orphaned \"*\" footnote marker that\n// used to point at the (now removed) corresponding author's superscript \"*\".\n// The rest of the \"Correspondence:\" line (the e-mail address) is left intact.\n\nconst body = context.document.body;\n\n// 1) Find every paragraph whose style is \"AuthorNames\" or \"Affiliations\" and\n//    delete it outright (these are the \"Firstname Lastname ...\" author line\n//    and the two \"<sup>n</sup> Affiliation n\" lines).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,style\");\nawait context.sync();\n\nconst paragraphsToDelete = [];\nlet correspondenceParagraph = null;\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.style === \"AuthorNames\" || paragraph.style === \"Affiliations\") {\n    paragraphsToDelete.push(paragraph);\n  } else if (paragraph.style === \"Correspondence\" && !correspondenceParagraph) {\n    correspondenceParagraph = paragraph;\n  }\n}\n\nfor (const paragraph of paragraphsToDelete) {\n  paragraph.delete();\n}\nawait context.sync();\n\n// 2) In the \"Correspondence\" paragraph, remove the leading bold \"*\" run that\n//    used to mark the corresponding author, leaving \" Correspondence: e-mail@e-mail.com\".\nif (correspondenceParagraph) {\n  const searchResults = correspondenceParagraph.search(\"*\", { matchWildcards: false });\n  searchResults.load(\"items,text,font\");\n  await context.sync();\n\n  for (const result of searchResults.items) {\n    if (result.text === \"*\") {\n      result.font.load(\"bold\");\n      await context.sync();\n      if (result.font.bold) {\n        result.delete();\n      }\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the placeholder \"AuthorNames\" and \"Affiliations\" paragraphs from the\n# document header block, and drop the now-orphaned \"*\" footnote marker that\n# used to point at the (now removed) corresponding author's superscript \"*\".\n# The rest of the \"Correspondence:\" line (the e-mail address) is left intact.\n\n$d = $word.ActiveDocument\n\n# 1) Collect every paragraph whose style is \"AuthorNames\" or \"Affiliations\"\n#    (the \"Firstname Lastname ...\" author line and the two\n#    \"<sup>n</sup> Affiliation n\" lines).\n$paragraphsToDelete = @()\nforeach ($p in $d.Paragraphs) {\n    $styleName = $p.Style.NameLocal\n    if ($styleName -eq \"AuthorNames\" -or $styleName -eq \"Affiliations\") {\n        $paragraphsToDelete += $p\n    }\n}\n\n# Delete them back-to-front: deleting earlier paragraphs first would shift\n# the positions backing the later Paragraph references in this array.\nfor ($i = $paragraphsToDelete.Count - 1; $i -ge 0; $i--) {\n    $paragraphsToDelete[$i].Range.Delete()\n}\n\n# 2) In the \"Correspondence\" paragraph, remove the leading bold \"*\" run that\n#    used to mark the corresponding author, leaving\n#    \" Correspondence: e-mail@e-mail.com\".\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"Correspondence\") {\n        $firstChar = $p.Range.Characters(1)\n        if ($firstChar.Text -eq \"*\" -and $firstChar.Font.Bold) {\n            $firstChar.Delete()\n        }\n        break\n    }\n}\n"}
